$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update CTC value which drives all the dependent tax formulas
$ws.Range("B2").Value = 1800000

# Update the active cell selection to match the saved view state
$ws.Range("G13").Select()
